# edit.ps1
# Applies the "New crime data collected" update to the 81st Precinct
# CompStat weekly report:
#   - Bumps the report header (Volume/Number + week-covering date range).
#   - Refreshes the Crime Complaints table (rows 14-26, 28-31) with the
#     newly collected weekly/28-day/YTD/2-year figures and their derived
#     percent-change columns.
#   - A handful of cells that were previously "not yet available"
#     placeholders (displayed as "0"/"***.*") now have real figures, so
#     their number format is (re)applied before the value is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header -------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/9/2024  Through  9/15/2024"

# --- Crime Complaints table -----------------------------------------
    # Row 14
    $ws.Range("M14").Value = -71.428571428571
    # Row 15
    $ws.Range("M15").Value = -50
    $ws.Range("N15").Value = -84.482758620689
    # Row 16
    $ws.Range("C16").Value = 5
    $ws.Range("E16").Value = 400
    $ws.Range("G16").Value = 4
    $ws.Range("H16").Value = 100
    $ws.Range("I16").Value = 90
    $ws.Range("J16").Value = 118
    $ws.Range("K16").Value = -23.728813559322
    $ws.Range("L16").Value = -25.619834710743
    $ws.Range("M16").Value = -61.538461538461
    $ws.Range("N16").Value = -89.221556886227
    # Row 17
    $ws.Range("D17").Value = 5
    $ws.Range("E17").Value = -60
    $ws.Range("G17").Value = 13
    $ws.Range("H17").Value = 38.461538461538
    $ws.Range("I17").Value = 216
    $ws.Range("J17").Value = 223
    $ws.Range("K17").Value = -3.139013452914
    $ws.Range("L17").Value = -11.475409836065
    $ws.Range("M17").Value = -8.860759493670
    $ws.Range("N17").Value = -66.769230769230
    # Row 18
    $ws.Range("D18").Value = 1
    $ws.Range("G18").Value = 11
    $ws.Range("H18").Value = -72.727272727272
    $ws.Range("J18").Value = 112
    $ws.Range("K18").Value = -50.892857142857
    $ws.Range("L18").Value = -60.144927536231
    $ws.Range("M18").Value = -64.516129032258
    $ws.Range("N18").Value = -90.107913669064
    # Row 19
    $ws.Range("C19").Value = 8
    $ws.Range("E19").Value = 33.333333333333
    $ws.Range("F19").Value = 16
    $ws.Range("G19").Value = 23
    $ws.Range("H19").Value = -30.434782608695
    $ws.Range("I19").Value = 200
    $ws.Range("J19").Value = 239
    $ws.Range("K19").Value = -16.317991631799
    $ws.Range("L19").Value = -32.432432432432
    $ws.Range("M19").Value = -11.111111111111
    $ws.Range("N19").Value = -24.242424242424
    # Row 20
    $ws.Range("C20").Value = 3
    $ws.Range("D20").Value = 4
    $ws.Range("E20").Value = -25
    $ws.Range("G20").Value = 13
    $ws.Range("H20").Value = -23.076923076923
    $ws.Range("I20").Value = 68
    $ws.Range("J20").Value = 78
    $ws.Range("K20").Value = -12.820512820512
    $ws.Range("L20").Value = -32.673267326732
    $ws.Range("M20").Value = 4.615384615384
    $ws.Range("N20").Value = -84.367816091954
    # Row 21
    $ws.Range("C21").Value = 18
    $ws.Range("E21").Value = 5.882352941176
    $ws.Range("F21").Value = 56
    $ws.Range("G21").Value = 65
    $ws.Range("H21").Value = -13.846153846153
    $ws.Range("I21").Value = 642
    $ws.Range("J21").Value = 790
    $ws.Range("K21").Value = -18.734177215189
    $ws.Range("L21").Value = -30.293159609120
    $ws.Range("M21").Value = -32.278481012658
    $ws.Range("N21").Value = -77.217885024840
    # Row 22
    $ws.Range("F22").Value = 2
    $ws.Range("I22").Value = 17
    $ws.Range("K22").Value = 70
    $ws.Range("L22").Value = 70
    $ws.Range("M22").Value = -10.526315789473
    # Row 23
    $ws.Range("C23").NumberFormat = "#,##0"
    $ws.Range("C23").Value = 1
    $ws.Range("D23").Value = 3
    $ws.Range("E23").Value = -66.666666666666
    $ws.Range("F23").Value = 2
    $ws.Range("G23").Value = 7
    $ws.Range("H23").Value = -71.428571428571
    $ws.Range("I23").Value = 66
    $ws.Range("J23").Value = 65
    $ws.Range("K23").Value = 1.538461538461
    $ws.Range("L23").Value = 6.451612903225
    $ws.Range("M23").Value = 11.864406779661
    # Row 24
    $ws.Range("C24").Value = 19
    $ws.Range("D24").Value = 26
    $ws.Range("E24").Value = -26.923076923076
    $ws.Range("F24").Value = 66
    $ws.Range("G24").Value = 79
    $ws.Range("H24").Value = -16.455696202531
    $ws.Range("I24").Value = 501
    $ws.Range("J24").Value = 619
    $ws.Range("K24").Value = -19.063004846526
    $ws.Range("L24").Value = -15.514333895446
    $ws.Range("M24").Value = -10.375670840787
    # Row 25
    $ws.Range("C25").Value = 3
    $ws.Range("D25").Value = 4
    $ws.Range("E25").Value = -25
    $ws.Range("F25").Value = 17
    $ws.Range("G25").Value = 11
    $ws.Range("H25").Value = 54.545454545454
    $ws.Range("I25").Value = 87
    $ws.Range("J25").Value = 158
    $ws.Range("K25").Value = -44.936708860759
    $ws.Range("L25").Value = -52.972972972973
    # Row 26
    $ws.Range("C26").Value = 10
    $ws.Range("D26").Value = 5
    $ws.Range("E26").Value = 100
    $ws.Range("F26").Value = 38
    $ws.Range("G26").Value = 28
    $ws.Range("H26").Value = 35.714285714285
    $ws.Range("I26").Value = 304
    $ws.Range("J26").Value = 358
    $ws.Range("K26").Value = -15.083798882681
    $ws.Range("L26").Value = 2.356902356902
    $ws.Range("M26").Value = -47.130434782608
    # Row 28
    $ws.Range("D28").Value = 4
    $ws.Range("G28").Value = 9
    $ws.Range("H28").Value = -88.888888888888
    $ws.Range("J28").Value = 24
    $ws.Range("K28").Value = 16.666666666666
    # Row 29
    $ws.Range("D29").NumberFormat = "#,##0"
    $ws.Range("D29").Value = 2
    $ws.Range("E29").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("E29").Value = -100
    $ws.Range("G29").Value = 2
    $ws.Range("H29").Value = -50
    $ws.Range("J29").Value = 13
    $ws.Range("K29").Value = 23.076923076923
    $ws.Range("M29").Value = -64.444444444444
    $ws.Range("N29").Value = -88.489208633093
    # Row 30
    $ws.Range("D30").NumberFormat = "#,##0"
    $ws.Range("D30").Value = 1
    $ws.Range("E30").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("E30").Value = -100
    $ws.Range("J30").Value = 10
    $ws.Range("K30").Value = 40
    $ws.Range("M30").Value = -63.157894736842
    $ws.Range("N30").Value = -88.617886178861
    # Row 31
    $ws.Range("D31").NumberFormat = "#,##0"
    $ws.Range("D31").Value = 1
    $ws.Range("E31").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("E31").Value = -100
    $ws.Range("G31").NumberFormat = "#,##0"
    $ws.Range("G31").Value = 1
    $ws.Range("H31").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("H31").Value = -100
    $ws.Range("J31").Value = 2
